$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns F:G use the same "text" number format as the existing data (column E's style)
$ws.Range("F1:G7").NumberFormat = "@"

# Add the new header cells in F1/G1
$ws.Range("F1").Value = "firstname"
$ws.Range("G1").Value = "lastname"

# Fill data rows 2..7 with "sindhu" / "boston"
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 6).Value = "sindhu"
    $ws.Cells.Item($r, 7).Value = "boston"
}

# Update the selection to match the post-edit state
$ws.Range("F3:G7").Select()
